# Add the new "2022-Q4" quarterly sheet right after "总计", and update the
# summary ("总计") sheet with the new quarter's totals.
#
# All other quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2,
# 2021-Q1, 2020-Q4) keep their own data untouched - inserting the new sheet
# simply shifts them one tab to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row right under the
#    header for the "2022-Q4" quarter, push every other row down by one,
#    and renumber the leading index column.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# carry the index-column style (bordered/bold/centered) down into the new row
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.29

for ($r = 3; $r -le 9; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert the brand-new "2022-Q4" worksheet right after "总计" and fill
#    in its fund-holdings table.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "159667"
$q4.Range("C2").Value = "国泰中证机床ETF"
$q4.Range("D2").Value = "3.49"
$q4.Range("E2").Value = "99.21"
$q4.Range("F2").Value = "4.21"
$q4.Range("G2").Value = "0.1469"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "015986"
$q4.Range("C3").Value = "中海新兴成长六个月持有期混合"
$q4.Range("D3").Value = "2.56"
$q4.Range("E3").Value = "40.87"
$q4.Range("F3").Value = "2.66"
$q4.Range("G3").Value = "0.0681"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "159663"
$q4.Range("C4").Value = "华夏中证机床ETF"
$q4.Range("D4").Value = "1.01"
$q4.Range("E4").Value = "97.42"
$q4.Range("F4").Value = "4.19"
$q4.Range("G4").Value = "0.0423"
$q4.Range("H4").Value = 8

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "000166"
$q4.Range("C5").Value = "中海信息产业精选混合"
$q4.Range("D5").Value = "0.72"
$q4.Range("E5").Value = "83.53"
$q4.Range("F5").Value = "3.87"
$q4.Range("G5").Value = "0.0279"
$q4.Range("H5").Value = 5
